$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 'last updated' timestamp shown in A1
$ws.Range("A1").Value = 'Datos actualizados a 23 de Marzo de 2020 a las 12:16'

# Country order + updated stats (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) for rows 4-195
$rows = @(
    @(4, 'China', 81093, 39, 72703, 5120, 1749, 9, 3270),
    @(5, 'Italia', 59138, 0, 7024, 46638, 3000, 0, 5476),
    @(6, 'Estados Unidos', 35070, 1524, 178, 34434, 795, 39, 458),
    @(7, 'España', 33089, 4321, 3355, 27552, 2355, 410, 2182),
    @(8, 'Alemania', 26220, 1347, 266, 25843, 23, 17, 111),
    @(9, 'Iran', 23049, 1411, 8376, 12861, 0, 127, 1812),
    @(10, 'Francia', 16018, 0, 2200, 13144, 1746, 0, 674),
    @(11, 'Corea del Sur', 8961, 64, 3166, 5684, 59, 7, 111),
    @(12, 'Suiza', 8234, 760, 131, 8001, 141, 4, 102),
    @(13, 'Reino Unido', 5683, 0, 135, 5267, 20, 0, 281),
    @(14, 'Paises Bajos', 4204, 0, 2, 4023, 354, 0, 179),
    @(15, 'Austria', 3759, 177, 9, 3734, 14, 0, 16),
    @(16, 'Belgica', 3743, 342, 350, 3305, 322, 13, 88),
    @(17, 'Noruega', 2415, 30, 6, 2401, 32, 1, 8),
    @(18, 'Suecia', 1934, 0, 16, 1893, 80, 4, 25),
    @(19, 'Australia', 1709, 100, 88, 1614, 11, 0, 7),
    @(20, 'Brasil', 1604, 58, 2, 1577, 18, 0, 25),
    @(21, 'Portugal', 1600, 0, 5, 1581, 26, 0, 14),
    @(22, 'Malasia', 1518, 212, 159, 1345, 57, 4, 14),
    @(23, 'Canada', 1470, 0, 14, 1436, 1, 0, 20),
    @(24, 'Dinamarca', 1450, 55, 1, 1436, 55, 0, 13),
    @(25, 'Israel', 1238, 167, 37, 1200, 24, 0, 1),
    @(26, 'Turquia', 1236, 0, 0, 1206, 0, 0, 30),
    @(27, 'Chequia', 1165, 45, 6, 1158, 19, 0, 1),
    @(28, 'Japon', 1101, 0, 235, 825, 49, 0, 41),
    @(29, 'Irlanda', 906, 0, 5, 897, 29, 0, 4),
    @(30, 'Pakistan', 804, 28, 13, 785, 0, 1, 6),
    @(31, 'Luxemburgo', 798, 0, 6, 784, 3, 0, 8),
    @(32, 'Ecuador', 789, 0, 3, 772, 2, 0, 14),
    @(33, 'Tailandia', 721, 122, 52, 668, 7, 0, 1),
    @(34, 'Crucero', 712, 0, 567, 137, 15, 0, 8),
    @(35, 'Finlandia', 686, 60, 10, 675, 12, 0, 1),
    @(36, 'Polonia', 649, 15, 13, 629, 3, 0, 7),
    @(37, 'Chile', 632, 0, 8, 623, 7, 0, 1),
    @(38, 'Grecia', 624, 0, 19, 589, 18, 1, 16),
    @(39, 'Indonesia', 579, 65, 30, 500, 0, 1, 49),
    @(40, 'Rumania', 576, 143, 73, 499, 14, 1, 4),
    @(41, 'Islandia', 568, 0, 5, 562, 1, 0, 1),
    @(42, 'Arabia Saudita', 511, 0, 17, 494, 0, 0, 0),
    @(43, 'Catar', 494, 0, 33, 461, 6, 0, 0),
    @(44, 'Filipinas', 462, 82, 18, 411, 1, 8, 33),
    @(45, 'Singapur', 455, 0, 144, 309, 14, 0, 2),
    @(46, 'Rusia', 438, 71, 16, 421, 0, 0, 1),
    @(47, 'India', 425, 29, 24, 393, 0, 1, 8),
    @(48, 'Eslovenia', 414, 0, 0, 411, 12, 1, 3),
    @(49, 'Peru', 363, 0, 1, 357, 5, 0, 5),
    @(50, 'Hong Kong', 356, 38, 100, 252, 4, 0, 4),
    @(51, 'Estonia', 352, 26, 4, 348, 4, 0, 0),
    @(52, 'Barein', 337, 3, 160, 175, 3, 0, 2),
    @(53, 'Egipto', 327, 0, 56, 257, 0, 0, 14),
    @(54, 'Mexico', 316, 65, 4, 310, 1, 0, 2),
    @(55, 'Panama', 313, 0, 1, 309, 7, 0, 3),
    @(56, 'Croacia', 306, 52, 5, 300, 5, 0, 1),
    @(57, 'Sudafrica', 274, 0, 2, 272, 0, 0, 0),
    @(58, 'Argentina', 266, 0, 27, 235, 0, 0, 4),
    @(59, 'Libano', 256, 8, 8, 244, 4, 0, 4),
    @(60, 'Colombia', 235, 4, 3, 230, 0, 0, 2),
    @(61, 'Irak', 233, 0, 57, 156, 0, 0, 20),
    @(62, 'Serbia', 222, 0, 2, 218, 4, 0, 2),
    @(63, 'Republica Dominicana', 202, 0, 0, 199, 0, 0, 3),
    @(64, 'Argelia', 201, 0, 65, 119, 0, 0, 17),
    @(65, 'Taiwan', 195, 26, 28, 165, 0, 0, 2),
    @(66, 'Armenia', 194, 0, 2, 192, 6, 0, 0),
    @(67, 'Bulgaria', 190, 3, 3, 184, 3, 0, 3),
    @(68, 'Kuwait', 189, 1, 30, 159, 5, 0, 0),
    @(69, 'Eslovaquia', 186, 1, 7, 179, 2, 0, 0),
    @(70, 'Letonia', 180, 41, 1, 179, 0, 0, 0),
    @(71, 'San Marino', 175, 0, 4, 151, 13, 0, 20),
    @(72, 'Hungria', 167, 36, 16, 144, 6, 1, 7),
    @(73, 'Uruguay', 158, 0, 0, 158, 2, 0, 0),
    @(74, 'Lituania', 154, 11, 1, 152, 1, 0, 1),
    @(75, 'Emiratos Arabes Unidos', 153, 0, 38, 113, 2, 0, 2),
    @(76, 'Costa Rica', 134, 0, 2, 130, 2, 0, 2),
    @(77, 'Bosnia y Herzegovina', 128, 2, 2, 125, 1, 0, 1),
    @(78, 'Marruecos', 122, 7, 3, 115, 1, 0, 4),
    @(79, 'Vietnam', 121, 8, 17, 104, 2, 0, 0),
    @(80, 'Islas Feroe', 118, 3, 14, 104, 0, 0, 0),
    @(81, 'Republica de Macedonia', 115, 0, 1, 112, 1, 1, 2),
    @(82, 'Principado de Andorra', 113, 0, 1, 111, 2, 0, 1),
    @(83, 'Jordania', 112, 0, 1, 111, 0, 0, 0),
    @(84, 'Nueva Zelanda', 102, 0, 0, 102, 0, 0, 0),
    @(85, 'Republica de Chipre', 95, 0, 3, 91, 3, 0, 1),
    @(86, 'Moldavia', 94, 0, 2, 91, 3, 0, 1),
    @(87, 'Brunei', 91, 3, 2, 89, 2, 0, 0),
    @(88, 'Malta', 90, 0, 2, 88, 1, 0, 0),
    @(89, 'Albania', 89, 0, 2, 82, 2, 3, 5),
    @(90, 'Sri Lanka', 87, 5, 3, 84, 2, 0, 0),
    @(91, 'Camboya', 86, 2, 2, 84, 1, 0, 0),
    @(92, 'Bielorrusia', 81, 5, 22, 59, 0, 0, 0),
    @(93, 'Venezuela', 77, 7, 15, 62, 2, 0, 0),
    @(94, 'Tunez', 75, 0, 1, 71, 7, 0, 3),
    @(95, 'Burkina Faso', 75, 0, 5, 66, 0, 0, 4),
    @(96, 'Ucrania', 73, 0, 1, 69, 0, 0, 3),
    @(97, 'Senegal', 67, 0, 5, 62, 0, 0, 0),
    @(98, 'Oman', 66, 11, 17, 49, 0, 0, 0),
    @(99, 'Azerbaiyan', 65, 0, 11, 53, 0, 0, 1),
    @(100, 'Reunion', 64, 0, 1, 63, 0, 0, 0),
    @(101, 'Kazajistan', 62, 2, 0, 62, 0, 0, 0),
    @(102, 'Estado de Palestina', 59, 0, 17, 42, 0, 0, 0),
    @(103, 'Guadalupe', 58, 0, 0, 57, 4, 0, 1),
    @(104, 'Camerun', 56, 16, 2, 54, 0, 0, 0),
    @(105, 'Georgia', 54, 0, 8, 46, 1, 0, 0),
    @(106, 'Trinidad yTobago', 50, 0, 0, 50, 0, 0, 0),
    @(107, 'Liechtenstein', 46, 9, 0, 46, 0, 0, 0),
    @(108, 'Uzbekistan', 46, 3, 0, 46, 0, 0, 0),
    @(109, 'Martinica', 44, 0, 0, 43, 7, 0, 1),
    @(110, 'Afganistan', 40, 0, 1, 38, 0, 0, 1),
    @(111, 'Nigeria', 36, 6, 2, 33, 0, 1, 1),
    @(112, 'Cuba', 35, 0, 0, 34, 0, 0, 1),
    @(113, 'Banglades', 33, 6, 3, 27, 0, 1, 3),
    @(114, 'Consejo Danes para los Refugiados', 30, 0, 0, 29, 0, 0, 1),
    @(115, 'Guam', 29, 2, 0, 28, 0, 0, 1),
    @(116, 'Mauricio', 28, 0, 0, 26, 1, 0, 2),
    @(117, 'Bolivia', 27, 3, 0, 27, 0, 0, 0),
    @(118, 'Honduras', 26, 0, 0, 26, 0, 0, 0),
    @(119, 'Costa de Marfil', 25, 11, 2, 23, 0, 0, 0),
    @(120, 'Ghana', 24, 1, 0, 23, 0, 0, 1),
    @(121, 'Macao', 24, 2, 10, 14, 0, 0, 0),
    @(122, 'Monaco', 23, 0, 1, 22, 0, 0, 0),
    @(123, 'Puerto Rico', 23, 0, 0, 21, 0, 1, 2),
    @(124, 'Paraguay', 22, 0, 0, 21, 1, 0, 1),
    @(125, 'Montenegro', 22, 1, 0, 21, 0, 1, 1),
    @(126, 'Mayotte', 21, 10, 0, 21, 0, 0, 0),
    @(127, 'Ruanda', 19, 0, 0, 19, 0, 0, 0),
    @(128, 'Guyana', 19, 0, 0, 18, 0, 0, 1),
    @(129, 'Guatemala', 19, 0, 0, 18, 0, 0, 1),
    @(130, 'Jamaica', 19, 0, 2, 16, 0, 0, 1),
    @(131, 'Guayana Francesa', 18, 0, 0, 18, 0, 0, 0),
    @(132, 'Polinesia Francesa', 18, 0, 0, 18, 0, 0, 0),
    @(133, 'Barbados', 17, 3, 0, 17, 0, 0, 0),
    @(134, 'Togo', 16, 0, 0, 16, 0, 0, 0),
    @(135, 'Kenia', 15, 0, 0, 15, 0, 0, 0),
    @(136, 'Gibraltar', 15, 0, 2, 13, 0, 0, 0),
    @(137, 'Kirguistan', 14, 0, 0, 14, 0, 0, 0),
    @(138, 'Maldivas', 13, 0, 5, 8, 0, 0, 0),
    @(139, 'Tanzania', 12, 0, 0, 12, 0, 0, 0),
    @(140, 'Madagascar', 12, 9, 0, 12, 0, 0, 0),
    @(141, 'Etiopia', 11, 0, 0, 11, 0, 0, 0),
    @(142, 'Mongolia', 10, 0, 0, 10, 0, 0, 0),
    @(143, 'Aruba', 9, 0, 1, 8, 0, 0, 0),
    @(144, 'Nueva Caledonia', 8, 4, 0, 8, 0, 0, 0),
    @(145, 'Seychelles', 7, 0, 0, 7, 0, 0, 0),
    @(146, 'Islas Virgenes de los Estados Unidos', 6, 0, 0, 6, 0, 0, 0),
    @(147, 'Bermudas', 6, 0, 0, 6, 0, 0, 0),
    @(148, 'Guinea Ecuatorial', 6, 0, 0, 6, 0, 0, 0),
    @(149, 'Isla de Man', 5, 0, 0, 5, 0, 0, 0),
    @(150, 'Surinam', 5, 0, 0, 5, 0, 0, 0),
    @(151, 'Haiti', 5, 3, 0, 5, 0, 0, 0),
    @(152, 'San Martin (Parte Francesa)', 5, 0, 0, 5, 0, 0, 0),
    @(153, 'Gabon', 5, 0, 0, 4, 0, 0, 1),
    @(154, 'Groenlandia', 4, 0, 0, 4, 0, 0, 0),
    @(155, 'Suazilandia', 4, 0, 0, 4, 0, 0, 0),
    @(156, 'Bahamas', 4, 0, 0, 4, 0, 0, 0),
    @(157, 'Guinea', 4, 2, 0, 4, 0, 0, 0),
    @(158, 'El Salvador', 3, 0, 0, 3, 0, 0, 0),
    @(159, 'Zimbabue', 3, 0, 0, 3, 0, 0, 0),
    @(160, 'Republica de Africa Central', 3, 0, 0, 3, 0, 0, 0),
    @(161, 'Zambia', 3, 0, 0, 3, 0, 0, 0),
    @(162, 'Liberia', 3, 0, 0, 3, 0, 0, 0),
    @(163, 'Congo', 3, 0, 0, 3, 0, 0, 0),
    @(164, 'Fiyi', 3, 1, 0, 3, 0, 0, 0),
    @(165, 'San Bartolome', 3, 0, 0, 3, 0, 0, 0),
    @(166, 'Namibia', 3, 0, 0, 3, 0, 0, 0),
    @(167, 'Cabo Verde', 3, 0, 0, 3, 0, 0, 0),
    @(168, 'Curazao', 3, 0, 0, 2, 0, 0, 1),
    @(169, 'Islas Caimanes', 3, 0, 0, 2, 0, 0, 1),
    @(170, 'Butan', 2, 0, 0, 2, 0, 0, 0),
    @(171, 'Nicaragua', 2, 0, 0, 2, 0, 0, 0),
    @(172, 'Mauritania', 2, 0, 0, 2, 0, 0, 0),
    @(173, 'Benin', 2, 0, 0, 2, 0, 0, 0),
    @(174, 'Santa Lucia', 2, 0, 0, 2, 0, 0, 0),
    @(175, 'Angola', 2, 0, 0, 2, 0, 0, 0),
    @(176, 'Niger', 2, 0, 0, 2, 0, 0, 0),
    @(177, 'Sudan', 2, 0, 0, 1, 0, 0, 1),
    @(178, 'Nepal', 2, 1, 1, 1, 0, 0, 0),
    @(179, 'Uganda', 1, 0, 0, 1, 0, 0, 0),
    @(180, 'Mozambique', 1, 0, 0, 1, 0, 0, 0),
    @(181, 'Somalia', 1, 0, 0, 1, 0, 0, 0),
    @(182, 'Republica de Yibuti', 1, 0, 0, 1, 0, 0, 0),
    @(183, 'Siria', 1, 0, 0, 1, 0, 0, 0),
    @(184, 'San Vicente y las Granadinas', 1, 0, 0, 1, 0, 0, 0),
    @(185, 'Antigua y Barbuda', 1, 0, 0, 1, 0, 0, 0),
    @(186, 'Republica del Chad', 1, 0, 0, 1, 0, 0, 0),
    @(187, 'Montserrat', 1, 0, 0, 1, 0, 0, 0),
    @(188, 'Dominica', 1, 0, 0, 1, 0, 0, 0),
    @(189, 'Granada', 1, 0, 0, 1, 0, 0, 0),
    @(190, 'Santa Sede', 1, 0, 0, 1, 0, 0, 0),
    @(191, 'San Martin (Parte Holandesa)', 1, 0, 0, 1, 0, 0, 0),
    @(192, 'Papua Nueva Guinea', 1, 0, 0, 1, 0, 0, 0),
    @(193, 'Timor Oriental', 1, 0, 0, 1, 0, 0, 0),
    @(194, 'Eritrea', 1, 0, 0, 1, 0, 0, 0),
    @(195, 'Gambia', 1, 0, 0, 0, 0, 1, 1)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}

Write-Output "Updated $($rows.Count) rows"
